$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings are preserved exactly
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.619.36'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '3.328.40'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '580.50'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '175.32'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').Value = '3.325.11'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('D10').Value = '0.179'
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('D11').Value = '0.579'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '46.36'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '704.45'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = '3.882.12'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '8.44'
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('D17').Value = '67.668.36'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '3.336.79'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').Value = '17.35'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').Value = '10.97'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').Value = '5.41'
$ws.Range('E23').Value = '  +3.70%  '
$ws.Range('D24').Value = '16.90'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').Value = '98.18'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').Value = '2.68'
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('D28').Value = '9.47'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').Value = '33.19'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '8.53'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  +4.35%  '
$ws.Range('D32').Value = '570.38'
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('D33').Value = '10.97'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').Value = '57.23'
$ws.Range('E35').Value = '  +3.33%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '3.700.22'
$ws.Range('E37').Value = '  -4.48%  '
$ws.Range('D38').Value = '3.32'
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('D39').Value = '33.99'
$ws.Range('E39').Value = '  +4.98%  '
$ws.Range('D40').Value = '0.130'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = '3.19'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = '2.64'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '0.336'
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0672'
$ws.Range('E44').Value = '  -2.07%  '
$ws.Range('D45').Value = '3.29'
$ws.Range('E45').Value = '  -2.94%  '
$ws.Range('D46').Value = '0.0405'
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('E47').Value = '  +5.24%  '
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('E50').Value = '  -5.70%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '128.63'
$ws.Range('E51').Value = '  -0.15%  '
